$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.2719521585403843
$ws.Range("B2").Value = -0.9480512082778079
$ws.Range("A3").Value = -0.5046402318966351
$ws.Range("B3").Value = -1.328230418635144
$ws.Range("A4").Value = -0.4968758457351299
$ws.Range("B4").Value = -1.56067230255249
$ws.Range("A5").Value = -0.1343710193016882
$ws.Range("B5").Value = -0.4358971670648874
$ws.Range("A6").Value = -0.1901677104056745
$ws.Range("B6").Value = -0.3246907239788171
$ws.Range("A7").Value = -0.3246880953562151
$ws.Range("B7").Value = -1.047004261679337
$ws.Range("A8").Value = -0.5200882665034741
$ws.Range("B8").Value = -1.354288945185685
$ws.Range("A9").Value = -0.6658251095750938
$ws.Range("B9").Value = -1.820634497020904
$ws.Range("A10").Value = -0.3862821647651075
$ws.Range("B10").Value = -0.2122899687332847
$ws.Range("A11").Value = -0.2679831408532483
$ws.Range("B11").Value = -0.8366220560994736
$ws.Range("A12").Value = -0.1357295443320637
$ws.Range("B12").Value = -0.2459148758661844
$ws.Range("A13").Value = -0.04434679156126818
$ws.Range("B13").Value = -0.08690707323059295
$ws.Range("A14").Value = -0.2232668801970721
$ws.Range("B14").Value = -0.09222417103618752
$ws.Range("A15").Value = -0.06562417062016375
$ws.Range("B15").Value = 0.03951747570295799
$ws.Range("A16").Value = -0.2825648592767291
$ws.Range("B16").Value = -0.6617785069003083
$ws.Range("A17").Value = -0.05883505353179978
$ws.Range("B17").Value = 0.4180184793348665
$ws.Range("A18").Value = 0.1123348910098452
$ws.Range("B18").Value = 0.5876146098826521
$ws.Range("A19").Value = 0.1808141315001575
$ws.Range("B19").Value = 0.5210137682756716
$ws.Range("A20").Value = -0.2778179509908926
$ws.Range("B20").Value = -0.7169372435563657
$ws.Range("A21").Value = 0.06127891984245044
$ws.Range("B21").Value = 0.4595633979573926
$ws.Range("A22").Value = -0.08377492228719334
$ws.Range("B22").Value = 0.3300377871875838
$ws.Range("A23").Value = 0.08120439994254271
$ws.Range("B23").Value = 0.6023096494722882
$ws.Range("A24").Value = 0.8794172050174812
$ws.Range("B24").Value = 2.769622411896357
$ws.Range("A25").Value = 0.1895797039263419
$ws.Range("B25").Value = 0.7122379008115508
$ws.Range("A26").Value = 0.2137308596472731
$ws.Range("B26").Value = 0.6287290635470522
$ws.Range("A27").Value = 0.1449695361051999
$ws.Range("B27").Value = 0.6937888252873274
$ws.Range("A28").Value = 0.4230745641164471
$ws.Range("B28").Value = 1.354994384057687
$ws.Range("A29").Value = 0.6969601030214119
$ws.Range("B29").Value = 2.47392557490438
$ws.Range("A30").Value = 0.2384889439052985
$ws.Range("B30").Value = 0.8124472388371239
$ws.Range("A31").Value = 0.1412369496037052
$ws.Range("B31").Value = 0.6063764165710263
$ws.Range("A32").Value = 0.2035762874408452
$ws.Range("B32").Value = 1.075698624533763
$ws.Range("A33").Value = 0.08622761833012706
$ws.Range("B33").Value = 0.736304179896366
$ws.Range("A34").Value = 0.06523890904166356
$ws.Range("B34").Value = 0.4345939322966909
$ws.Range("A35").Value = 0.4704159447966034
$ws.Range("B35").Value = 1.314106427822339
$ws.Range("A36").Value = 0.2188211091452042
$ws.Range("B36").Value = 0.5236396228262162
$ws.Range("A37").Value = -0.0230313243254189
$ws.Range("B37").Value = -0.1470899638659299
$ws.Range("A38").Value = 0.2433689329518174
$ws.Range("B38").Value = 1.458571259270864
$ws.Range("A39").Value = -0.1139745991102089
$ws.Range("B39").Value = -0.3732865694291094
$ws.Range("A40").Value = 0.07844096572615572
$ws.Range("B40").Value = 0.578274097529577
$ws.Range("A41").Value = 0.0121108499720569
$ws.Range("B41").Value = 0.3731110870504676
$ws.Range("A42").Value = 0.3307377616317541
$ws.Range("B42").Value = 1.584342368571894
$ws.Range("A43").Value = -0.008975497306698567
$ws.Range("B43").Value = 0.5655516150151232
$ws.Range("A44").Value = 0.08766768056055341
$ws.Range("B44").Value = 0.3961495194008971
$ws.Range("A45").Value = -0.1190181052748567
$ws.Range("B45").Value = -0.07107604345940965
$ws.Range("A46").Value = -0.1761622947185779
$ws.Range("B46").Value = -0.4900795516462431
$ws.Range("A47").Value = -0.1787401865259877
$ws.Range("B47").Value = -0.3732203843117181
$ws.Range("A48").Value = -0.2274708112765429
$ws.Range("B48").Value = -0.5986664587414678
$ws.Range("A49").Value = -0.2261114141792094
$ws.Range("B49").Value = -0.6213793960768439
$ws.Range("A50").Value = -0.07579209810743265
$ws.Range("B50").Value = -0.1417118863480319
$ws.Range("A51").Value = -0.2516874376647317
$ws.Range("B51").Value = -0.6674725662163312
$ws.Range("A52").Value = -0.2516874376647317
$ws.Range("B52").Value = -0.6674725662163312
$ws.Range("A53").Value = -0.2116181978129459
$ws.Range("B53").Value = -0.450445904738707
$ws.Range("A54").Value = -0.2029489456450268
$ws.Range("B54").Value = -0.5347693233809983
$ws.Range("A55").Value = -0.1698827370189405
$ws.Range("B55").Value = -0.4159270119684224
$ws.Range("A56").Value = -0.09617086707026469
$ws.Range("B56").Value = -0.2617903839688462
$ws.Range("A57").Value = -0.1906883189435502
$ws.Range("B57").Value = -0.3661170676443674
$ws.Range("A58").Value = -0.1272461134363186
$ws.Range("B58").Value = -0.1836227512390941
$ws.Range("A59").Value = -0.2090118182953766
$ws.Range("B59").Value = -0.6763109878255052
$ws.Range("A60").Value = -0.2485642637880077
$ws.Range("B60").Value = -0.7153565954539691
$ws.Range("A61").Value = -0.2792493696662376
$ws.Range("B61").Value = -0.3578624880748587
$ws.Range("A62").Value = -0.1511779146553517
$ws.Range("B62").Value = -0.1379193774320158
$ws.Range("A63").Value = -0.514044960024056
$ws.Range("B63").Value = -1.297489951101472
$ws.Range("A64").Value = -0.319776463580188
$ws.Range("B64").Value = -0.7017671301135775
$ws.Range("A65").Value = -0.176894792088321
$ws.Range("B65").Value = -0.4869009484686553
$ws.Range("A66").Value = -0.07475894167698741
$ws.Range("B66").Value = -0.01139451757782755
$ws.Range("A67").Value = 0.07788704078412738
$ws.Range("B67").Value = 0.3707335271770831
